# Presentazione Progetto: bump the title placeholder's font size
# on slide 1 ("Progetto Ingegneria del Software / Gestione Ospedale")
# from 42pt to 54pt, matching the author's manual "increase font size"
# edit recorded in the change history (actId="255").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$titleShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Titolo 1") {
        $titleShape = $shp
        break
    }
}

# Select the whole title text and set its size to 54pt (was 42pt).
$titleShape.TextFrame.TextRange.Font.Size = 54
